# Apply updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: update D, E
$ws.Range("D2").Value = "25.960.26"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3: update D, E
$ws.Range("D3").Value = "1.638.10"
$ws.Range("E3").Value = "  -0.17%  "

# Row 4: update D, E
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.54%  "

# Row 5: update D, E
$ws.Range("D5").Value = "'214.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "

# Row 6: update D, E
$ws.Range("D6").Value = "'0.5102"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.18%  "

# Row 7: update D, E
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.68%  "

# Row 8: update D, E
$ws.Range("D8").Value = "'0.2559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.43%  "

# Row 9: update D, E
$ws.Range("D9").Value = "'0.06362"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "

# Row 10: update D, E
$ws.Range("D10").Value = "'19.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.20%  "

# Row 11: update D, E
$ws.Range("D11").Value = "'0.07754"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.79%  "

# Row 12: update D, E
$ws.Range("D12").Value = "'4.285"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "

# Row 13: update D, E
$ws.Range("D13").Value = "1.645.61"
$ws.Range("E13").Value = "  +0.00%  "

# Row 14: update D, E
$ws.Range("D14").Value = "'0.5433"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.06%  "

# Row 15: update D, E
$ws.Range("D15").Value = "0.0₅7756"
$ws.Range("E15").Value = "  -1.45%  "

# Row 16: update D, E
$ws.Range("D16").Value = "'64.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.45%  "

# Row 17: update D, E
$ws.Range("D17").Value = "25.962.73"
$ws.Range("E17").Value = "  +0.09%  "

# Row 18: update D, E
$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19: update D, E
$ws.Range("D19").Value = "'196.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.65%  "

# Row 20: update D, E
$ws.Range("D20").Value = "'4.425"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.16%  "

# Row 21: update D, E
$ws.Range("D21").Value = "'9.921"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.11%  "

# Row 22: update D, E
$ws.Range("D22").Value = "'6.037"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.24%  "

# Row 23: update D, E
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.79%  "

# Row 24: update D, E
$ws.Range("D24").Value = "'1.864"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.32%  "

# Row 25: update D, E
$ws.Range("D25").Value = "'141.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "

# Row 26: update D, E
$ws.Range("D26").Value = "'0.1194"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.60%  "

# Row 27: update D, E
$ws.Range("D27").Value = "'6.838"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.48%  "

# Row 28: update D, E
$ws.Range("D28").Value = "'15.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "

# Row 29: update D, E
$ws.Range("D29").Value = "'1.233"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.04%  "

# Row 30: update D, E
$ws.Range("D30").Value = "'0.04940"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "

# Row 31: update D, E
$ws.Range("D31").Value = "'3.243"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "

# Row 32: update D, E
$ws.Range("D32").Value = "'3.179"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "

# Row 33: update D, E
$ws.Range("D33").Value = "'1.528"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "

# Row 34: update D, E
$ws.Range("D34").Value = "'2.365"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.96%  "

# Row 35: update D, E
$ws.Range("D35").Value = "'0.8919"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "

# Row 36: update D, E
$ws.Range("D36").Value = "1.149.92"
$ws.Range("E36").Value = "  +1.69%  "

# Row 37: update D, E
$ws.Range("D37").Value = "'2.576"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.14%  "

# Row 38: update D, E
$ws.Range("D38").Value = "'0.5421"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.56%  "

# Row 39: update D, E
$ws.Range("D39").Value = "'0.01552"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.92%  "

# Row 40: update E
$ws.Range("E40").Value = "  -1.24%  "

# Row 41: update E
$ws.Range("E41").Value = "  -1.58%  "

# Row 42: update B, C, D, E
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8119"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.32%  "

# Row 43: update B, C, D, E
$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D43").Value = "0.0₈126"
$ws.Range("E43").Value = "  +7.67%  "

# Row 44: update B, C, D, E
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.457"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.97%  "

# Row 45: update B, C, D, E
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'99.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.57%  "

# Row 46: update D, E
$ws.Range("D46").Value = "1.776.81"
$ws.Range("E46").Value = "  +0.05%  "

# Row 47: update D, E
$ws.Range("D47").Value = "'0.4518"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.21%  "

# Row 48: update E
$ws.Range("E48").Value = "  -0.65%  "

# Row 49: update E
$ws.Range("E49").Value = "  -0.29%  "

# Row 50: update D, E
$ws.Range("D50").Value = "'0.05051"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.10%  "

# Row 51: update D, E
$ws.Range("D51").Value = "'1.001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.87%  "
